# customerTransactions: wipe out the sample order rows (2-8, columns A-H)
# but keep the order number that was already sitting in A2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:H8").ClearContents()
$ws.Range("A2").Value = 2
